$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typos in header row
$ws.Range("A1").Value = "NOM ETUDIANT"
$ws.Range("H1").Value = "DEVELOPEMENT WEB"

# Move active cell selection to A11
$ws.Range("A11").Select()
